$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a daily-updated time series (newest entries are inserted near
# the top, pushing older entries down). Two new price observations were
# added this edit:
#   - a new row inserted at row 42 (date 2022-03-18 / serial 44645)
#   - a new row inserted at row 100 (date 2022-03-17 / serial 44644)
# Everything else is an untouched, pure downward shift of the existing rows,
# which Rows.Insert() reproduces natively (formats + values move together).
# ---------------------------------------------------------------------------

# --- Insert new row #1 at sheet row 42 -------------------------------------
$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value  = 11
$ws.Cells.Item(42, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value  = "Bíobío"
$ws.Cells.Item(42, 4).Value  = 44645
$ws.Cells.Item(42, 5).Value  = 8
$ws.Cells.Item(42, 6).Value  = 100112043
$ws.Cells.Item(42, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(42, 8).Value  = "Sin especificar"
$ws.Cells.Item(42, 9).Value  = "Primera"
$ws.Cells.Item(42, 10).Value = 170
$ws.Cells.Item(42, 11).Value = 17000
$ws.Cells.Item(42, 12).Value = 18000
$ws.Cells.Item(42, 13).Value = 17529
$ws.Cells.Item(42, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(42, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(42, 16).Value = 292
$ws.Cells.Item(42, 17).Value = 60
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# --- Insert new row #2 at sheet row 100 (post first insert) ----------------
$ws.Rows.Item(100).Insert()

$ws.Cells.Item(100, 1).Value  = 11
$ws.Cells.Item(100, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(100, 3).Value  = "Bíobío"
$ws.Cells.Item(100, 4).Value  = 44644
$ws.Cells.Item(100, 5).Value  = 8
$ws.Cells.Item(100, 6).Value  = 100112043
$ws.Cells.Item(100, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(100, 8).Value  = "Sin especificar"
$ws.Cells.Item(100, 9).Value  = "Primera"
$ws.Cells.Item(100, 10).Value = 270
$ws.Cells.Item(100, 11).Value = 17000
$ws.Cells.Item(100, 12).Value = 18000
$ws.Cells.Item(100, 13).Value = 17444
$ws.Cells.Item(100, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(100, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value = 291
$ws.Cells.Item(100, 17).Value = 60
$ws.Cells.Item(100, 18).Value = "Hortaliza"

"done"
